$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 38, shifting the old rows 38-43 down to 40-45.
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()

# Fill the new row 38 with the new weekly data (Primera).
$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value = "Ñuble"
$ws.Cells.Item(38, 4).Value = 44617
$ws.Cells.Item(38, 5).Value = 16
$ws.Cells.Item(38, 6).Value = "Fruta"
$ws.Cells.Item(38, 7).Value = 100103
$ws.Cells.Item(38, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(38, 9).Value = 100103002
$ws.Cells.Item(38, 10).Value = "Ciruela"
$ws.Cells.Item(38, 11).Value = "Black Amber"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 60
$ws.Cells.Item(38, 14).Value = 10000
$ws.Cells.Item(38, 15).Value = 10000
$ws.Cells.Item(38, 16).Value = 10000
$ws.Cells.Item(38, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(38, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(38, 19).Value = 556
$ws.Cells.Item(38, 20).Value = 18

# Fill the new row 39 with the new weekly data (Segunda).
$ws.Cells.Item(39, 1).Value = 7
$ws.Cells.Item(39, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(39, 3).Value = "Ñuble"
$ws.Cells.Item(39, 4).Value = 44617
$ws.Cells.Item(39, 5).Value = 16
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100103
$ws.Cells.Item(39, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(39, 9).Value = 100103002
$ws.Cells.Item(39, 10).Value = "Ciruela"
$ws.Cells.Item(39, 11).Value = "Black Amber"
$ws.Cells.Item(39, 12).Value = "Segunda"
$ws.Cells.Item(39, 13).Value = 80
$ws.Cells.Item(39, 14).Value = 8000
$ws.Cells.Item(39, 15).Value = 9000
$ws.Cells.Item(39, 16).Value = 8500
$ws.Cells.Item(39, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(39, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(39, 19).Value = 472
$ws.Cells.Item(39, 20).Value = 18
